$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.494545936584473
$ws.Range("B1").Value = 1.567805171012878
$ws.Range("C1").Value = 3.717894315719604
$ws.Range("D1").Value = 2.32136344909668
$ws.Range("E1").Value = 0.8399289846420288
